$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$n = 72
$data = New-Object 'object[,]' $n,2
$data[0,0] = "Admiraal De Ruyter Ziekenhuis"
$data[0,1] = "definitief en/of vastgesteld RvB"
$data[1,0] = "Albert Schweitzer ziekenhuis"
$data[1,1] = "definitief en/of vastgesteld RvB"
$data[2,0] = "Alrijne Zorggroep"
$data[2,1] = "voorlopig"
$data[3,0] = "Amphia Ziekenhuis"
$data[3,1] = "voorlopig"
$data[4,0] = "Amsterdam UMC"
$data[4,1] = "definitief en/of vastgesteld RvB"
$data[5,0] = "Antoni van Leeuwenhoek"
$data[5,1] = "voorlopig"
$data[6,0] = "Antonius Zorggroep"
$data[6,1] = "definitief en/of vastgesteld RvB"
$data[7,0] = "Bravis ziekenhuis"
$data[7,1] = "definitief en/of vastgesteld RvB"
$data[8,0] = "Canisius-Wilhelmina Ziekenhuis"
$data[8,1] = "definitief en/of vastgesteld RvB"
$data[9,0] = "Catharina Ziekenhuis"
$data[9,1] = "definitief en/of vastgesteld RvB"
$data[10,0] = "CuraMare"
$data[10,1] = "voorlopig"
$data[11,0] = "Deventer Ziekenhuis"
$data[11,1] = "voorlopig"
$data[12,0] = "Diakonessenhuis"
$data[12,1] = "definitief en/of vastgesteld RvB"
$data[13,0] = "Dialysecentrum Ravenstein"
$data[13,1] = "voorlopig"
$data[14,0] = "Dianet"
$data[14,1] = "voorlopig"
$data[15,0] = "Elisabeth-TweeSteden Ziekenhuis"
$data[15,1] = "definitief en/of vastgesteld RvB"
$data[16,0] = "Elkerliek Ziekenhuis"
$data[16,1] = "definitief en/of vastgesteld RvB"
$data[17,0] = "Erasmus MC"
$data[17,1] = "definitief en/of vastgesteld RvB"
$data[18,0] = "Flevoziekenhuis"
$data[18,1] = "definitief en/of vastgesteld RvB"
$data[19,0] = "Franciscus Gasthuis & Vlietland Groep"
$data[19,1] = "definitief en/of vastgesteld RvB"
$data[20,0] = "Gelre Ziekenhuizen"
$data[20,1] = "definitief en/of vastgesteld RvB"
$data[21,0] = "Groene Hart Ziekenhuis"
$data[21,1] = "definitief en/of vastgesteld RvB"
$data[22,0] = "IJsselland Ziekenhuis"
$data[22,1] = "definitief en/of vastgesteld RvB"
$data[23,0] = "Isala Klinieken"
$data[23,1] = "definitief en/of vastgesteld RvB"
$data[24,0] = "Jeroen Bosch Ziekenhuis"
$data[24,1] = "definitief en/of vastgesteld RvB"
$data[25,0] = "LUMC"
$data[25,1] = "definitief en/of vastgesteld RvB"
$data[26,0] = "LangeLand Ziekenhuis"
$data[26,1] = "definitief en/of vastgesteld RvB"
$data[27,0] = "Laurentius Ziekenhuis"
$data[27,1] = "definitief en/of vastgesteld RvB"
$data[28,0] = "Maasstad ziekenhuis"
$data[28,1] = "definitief en/of vastgesteld RvB"
$data[29,0] = "Maastricht UMC+"
$data[29,1] = "definitief en/of vastgesteld RvB"
$data[30,0] = "Martini Ziekenhuis"
$data[30,1] = "definitief en/of vastgesteld RvB"
$data[31,0] = "Meander Medisch Centrum"
$data[31,1] = "definitief en/of vastgesteld RvB"
$data[32,0] = "Medisch Centrum Leeuwarden"
$data[32,1] = "definitief en/of vastgesteld RvB"
$data[33,0] = "Medisch Spectrum Twente"
$data[33,1] = "definitief en/of vastgesteld RvB"
$data[34,0] = "Máxima Medisch Centrum"
$data[34,1] = "definitief en/of vastgesteld RvB"
$data[35,0] = "Nij Smellinghe"
$data[35,1] = "definitief en/of vastgesteld RvB"
$data[36,0] = "Noordwest Ziekenhuisgroep"
$data[36,1] = "definitief en/of vastgesteld RvB"
$data[37,0] = "Ommelander Ziekenhuis Groep"
$data[37,1] = "definitief en/of vastgesteld RvB"
$data[38,0] = "Onze Lieve Vrouwe Gasthuis"
$data[38,1] = "definitief en/of vastgesteld RvB"
$data[39,0] = "Pantein"
$data[39,1] = "definitief en/of vastgesteld RvB"
$data[40,0] = "Prinses Máxima Centrum"
$data[40,1] = "definitief en/of vastgesteld RvB"
$data[41,0] = "Protestants Christelijk Ziekenhuis Ikazia"
$data[41,1] = "definitief en/of vastgesteld RvB"
$data[42,0] = "Radboudumc"
$data[42,1] = "definitief en/of vastgesteld RvB"
$data[43,0] = "Revalidatiecentrum Roessingh"
$data[43,1] = "voorlopig"
$data[44,0] = "Rivas Zorggroep"
$data[44,1] = "definitief en/of vastgesteld RvB"
$data[45,0] = "Rode Kruis Ziekenhuis"
$data[45,1] = "voorlopig"
$data[46,0] = "Slingeland Ziekenhuis (Santiz)"
$data[46,1] = "definitief en/of vastgesteld RvB"
$data[47,0] = "Spaarne Gasthuis"
$data[47,1] = "definitief en/of vastgesteld RvB"
$data[48,0] = "Spijkenisse Medisch Centrum"
$data[48,1] = "definitief en/of vastgesteld RvB"
$data[49,0] = "St Jansdal"
$data[49,1] = "definitief en/of vastgesteld RvB"
$data[50,0] = "St. Antonius Ziekenhuis"
$data[50,1] = "definitief en/of vastgesteld RvB"
$data[51,0] = "St. Jans Gasthuis"
$data[51,1] = "definitief en/of vastgesteld RvB"
$data[52,0] = "Stichting BovenIJ ziekenhuis"
$data[52,1] = "definitief en/of vastgesteld RvB"
$data[53,0] = "Stichting Dijklander Ziekenhuis"
$data[53,1] = "definitief en/of vastgesteld RvB"
$data[54,0] = "Stichting Reinier Haga Groep"
$data[54,1] = "definitief en/of vastgesteld RvB"
$data[55,0] = "Stichting Reinier de Graaf Groep"
$data[55,1] = "definitief en/of vastgesteld RvB"
$data[56,0] = "Stichting Rijnstate Ziekenhuis"
$data[56,1] = "definitief en/of vastgesteld RvB"
$data[57,0] = "Stichting St. Anna Zorggroep"
$data[57,1] = "definitief en/of vastgesteld RvB"
$data[58,0] = "Stichting Tergooi"
$data[58,1] = "definitief en/of vastgesteld RvB"
$data[59,0] = "Stichting VieCuri-Vitaal"
$data[59,1] = "definitief en/of vastgesteld RvB"
$data[60,0] = "Stichting Ziekenhuisgroep Twente (ZGT)"
$data[60,1] = "voorlopig"
$data[61,0] = "Streekziekenhuis Koningin Beatrix (Santiz)"
$data[61,1] = "definitief en/of vastgesteld RvB"
$data[62,0] = "Tjongerschans"
$data[62,1] = "definitief en/of vastgesteld RvB"
$data[63,0] = "Treant Zorggroep"
$data[63,1] = "voorlopig"
$data[64,0] = "UMC Utrecht"
$data[64,1] = "definitief en/of vastgesteld RvB"
$data[65,0] = "UMCG"
$data[65,1] = "definitief en/of vastgesteld RvB"
$data[66,0] = "Wilhelmina Ziekenhuis Assen"
$data[66,1] = "definitief en/of vastgesteld RvB"
$data[67,0] = "Ziekenhuis Amstelland"
$data[67,1] = "voorlopig"
$data[68,0] = "Ziekenhuis Bernhoven"
$data[68,1] = "definitief en/of vastgesteld RvB"
$data[69,0] = "Ziekenhuis De Gelderse Vallei"
$data[69,1] = "definitief en/of vastgesteld RvB"
$data[70,0] = "Ziekenhuis Rivierenland"
$data[70,1] = "definitief en/of vastgesteld RvB"
$data[71,0] = "Zuyderland Ziekenhuis"
$data[71,1] = "definitief en/of vastgesteld RvB"

$ws.Range("A2:B73").Value = $data

Write-Output "done"